# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Updates Price (D) / Volume(1h) (E) figures for each coin row, and re-applies
# the two rank-adjacent row swaps (ARBITRUM/Stellar and ordi/BitcoinSV) whose
# relative order flipped with this refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.667.41"
$ws.Range("E2").Value = "  +0.34%  "

$ws.Range("D3").Value = "2.474.97"
$ws.Range("E3").Value = "  -0.36%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.48"
$ws.Range("E5").Value = "  +1.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.60"
$ws.Range("E6").Value = "  -1.24%  "

$ws.Range("E7").Value = "  +1.93%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  +1.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0898"
$ws.Range("E10").Value = "  +14.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "32.94"
$ws.Range("E11").Value = "  +0.52%  "

$ws.Range("E12").Value = "  +0.36%  "

$ws.Range("D13").Value = "2.855.18"
$ws.Range("E13").Value = "  -0.38%  "

$ws.Range("E14").Value = "  +0.68%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.67"
$ws.Range("E15").Value = "  -2.57%  "

$ws.Range("D16").Value = "2.494.71"
$ws.Range("E16").Value = "  +2.87%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.791"
$ws.Range("E17").Value = "  +3.54%  "

$ws.Range("D18").Value = "41.608.63"
$ws.Range("E18").Value = "  +0.20%  "

$ws.Range("E19").Value = "  +2.85%  "

$ws.Range("E20").Value = "  +0.89%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.43"
$ws.Range("E21").Value = "  -0.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.45"
$ws.Range("E22").Value = "  +0.68%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.51"
$ws.Range("E23").Value = "  +1.84%  "

$ws.Range("E24").Value = "  +0.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.93"
$ws.Range("E25").Value = "  +1.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.87"
$ws.Range("E27").Value = "  -0.31%  "

$ws.Range("E28").Value = "  +4.45%  "

$ws.Range("E29").Value = "  +1.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.27"
$ws.Range("E30").Value = "  +0.60%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.59"
$ws.Range("E31").Value = "  -0.92%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.53"
$ws.Range("E32").Value = "  +0.42%  "

$ws.Range("E34").Value = "  +1.77%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.57"
$ws.Range("E35").Value = "  -0.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.43"
$ws.Range("E36").Value = "  -0.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.92"
$ws.Range("E37").Value = "  -0.42%  "

$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.83"
$ws.Range("E38").Value = "  -0.27%  "

$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.116"
$ws.Range("E39").Value = "  +1.25%  "

$ws.Range("E40").Value = "  -1.32%  "

$ws.Range("E41").Value = "  +1.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.98"
$ws.Range("E42").Value = "  -3.45%  "

$ws.Range("D43").Value = "1.980.20"
$ws.Range("E43").Value = "  -0.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.13"
$ws.Range("E44").Value = "  -3.28%  "

$ws.Range("E45").Value = "  +0.19%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.99"
$ws.Range("E46").Value = "  +0.88%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.15"
$ws.Range("E47").Value = "  +0.16%  "

$ws.Range("D48").Value = "2.708.71"
$ws.Range("E48").Value = "  -0.52%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.38"
$ws.Range("E49").Value = "  -0.83%  "

$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.47"
$ws.Range("E50").Value = "  -1.08%  "

$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.84"
$ws.Range("E51").Value = "  +1.88%  "
